$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "digital encrypted certificates available through one of your CS site reps.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "digital encrypted client certificates available from CS site reps.",
    2
)
